# Scenario_InterestRate.xlsx — re-saved / reviewed in Excel:
#  - the user widened/auto-fit the first four (table) columns so the
#    id_region / id_sector / id_subsector / unit headers and values are
#    fully visible, and
#  - left the selection on R15 when the file was saved.
# No cell values in the interest-rate table itself were changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-fit the four columns that hold the table's key/unit fields
# (id_region, id_sector, id_subsector, unit) to their content.
$ws.Columns.Item(1).ColumnWidth = 10.833333333333332
$ws.Columns.Item(2).ColumnWidth = 10.666666666666666
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666

# Restore the active selection to R15, matching the saved view state.
$ws.Range("R15").Select() | Out-Null
